$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.996.79'
$ws.Range('E2').Value = '  -2.05%  '
$ws.Range('D3').Value = '2.430.05'
$ws.Range('E3').Value = '  -0.57%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.55'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  -1.97%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.06'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  -2.74%  '
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.528'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  -0.67%  '
$ws.Range('D9').Value = '2.418.56'
$ws.Range('E9').Value = '  -0.96%  '
$ws.Range('E10').Value = '  +0.85%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.07'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  -2.75%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.338'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').Value = '  -2.09%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.07'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  -1.57%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000170'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  -1.72%  '
$ws.Range('D16').Value = '2.809.20'
$ws.Range('E16').Value = '  -2.51%  '
$ws.Range('D17').Value = '61.028.26'
$ws.Range('E17').Value = '  -1.80%  '
$ws.Range('D18').Value = '2.386.59'
$ws.Range('E18').Value = '  -2.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.55'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  -3.35%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.27'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  +1.78%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '322.21'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  -2.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.02'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  -2.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.08'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  +1.48%  '
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.87'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  -5.76%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '64.87'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  -1.47%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.71'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  -7.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '574.81'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  -7.24%  '
$ws.Range('D29').Value = '2.564.38'
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('D30').Value = '0.0₃0909'
$ws.Range('E30').Value = '  -5.06%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.85'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  -2.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.34'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  -6.54%  '
$ws.Range('E33').Value = '  -2.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.132'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  -6.71%  '
$ws.Range('E35').Value = '  +0.19%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.59'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  -6.75%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.367'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  -3.04%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.38'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  -4.26%  '
$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '148.67'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  -1.82%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.18'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  -0.84%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.06'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  -3.95%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '41.69'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  -1.67%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.65'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  -5.96%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.33'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  -6.08%  '
$ws.Range('D46').Value = '0.0₆0281'
$ws.Range('E46').Value = '  +18.33%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '140.71'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  -1.90%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.50'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  -3.78%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.593'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  -1.11%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0505'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  -4.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.34'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  -1.18%  '
